$d = $word.ActiveDocument

$replacements = @(
    @("339÷9=", "115÷7="),
    @("766÷3=", "931÷4="),
    @("355÷5=", "613÷5="),
    @("713÷9=", "580÷6="),
    @("399÷7=", "230÷6="),
    @("927÷2=", "628÷9="),
    @("354÷8=", "886÷5="),
    @("301÷7=", "670÷6="),
    @("105÷3=", "373÷7="),
    @("813÷7=", "260÷3="),
    @("157÷7=", "633÷8="),
    @("606÷2=", "855÷6="),
    @("134÷7=", "118÷8="),
    @("141÷6=", "562÷7="),
    @("453÷7=", "211÷2="),
    @("623÷6=", "960÷9="),
    @("759÷9=", "445÷5="),
    @("661÷6=", "393÷4="),
    @("887÷9=", "166÷5="),
    @("107÷4=", "898÷9="),
    @("289÷7=", "816÷8="),
    @("881÷3=", "196÷7="),
    @("295÷8=", "898÷7="),
    @("620÷9=", "985÷7="),
    @("113÷2=", "888÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
